$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Solar (column E) values for 2022 and 2024 rows with upstream-refreshed figures
$ws.Range("E24").Value = 144.58
$ws.Range("E26").Value = 215.87

$wb.Save()
